$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column O, shifting O:T to P:U
$ws.Range("O1:O2").EntireColumn.Insert()

# New column O header/value
$ws.Range("O1").Value = "Group.1"
$ws.Range("O2").Value = 6750225000

# Update recalculated values in B2:N2
$ws.Range("B2").Value = 10801998.33475939
$ws.Range("C2").Value = 2372282.164526457
$ws.Range("D2").Value = 26557345.39449133
$ws.Range("E2").Value = 1152364.165444766
$ws.Range("F2").Value = 8085423.201825836
$ws.Range("G2").Value = 1840377.530403443
$ws.Range("H2").Value = 2146872.416877256
$ws.Range("I2").Value = 10801998.33475939
$ws.Range("J2").Value = 45590812
$ws.Range("K2").Value = 123
$ws.Range("L2").Value = 28929627.55901779
$ws.Range("M2").Value = 9237787.367270602
$ws.Range("N2").Value = 3987249.947280699

# Updated values in the shifted P:U columns (formerly O:T)
$ws.Range("P2").Value = 42014.73273041348
$ws.Range("Q2").Value = 205086.7601341055
$ws.Range("R2").Value = 247101.492864519
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 109763.5155576515
$ws.Range("U2").Value = 109763.5155576515
